# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "59.226.82"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "2.998.88"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.50"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.25"
$ws.Range("E6").Value = "  +3.17%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "2.994.55"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("E11").Value = "  +7.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  +3.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.70"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.33"
$ws.Range("E16").Value = "  +7.34%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.494.89"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("D18").Value = "3.001.81"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "59.315.44"
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.03"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").Value = "  +5.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.13"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.73"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("B27").Value = "FirstDigitalUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  +10.58%  "
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.86"
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.81"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0997"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.97"
$ws.Range("E34").Value = "  +5.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("E35").Value = "  +5.27%  "
$ws.Range("D36").Value = "0.0₃0762"
$ws.Range("E36").Value = "  +9.95%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.90"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.64"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  +6.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "403.44"
$ws.Range("E41").Value = "  +5.42%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").Value = "2.767.39"
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.74"
$ws.Range("E46").Value = "  +25.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.86"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.53"
$ws.Range("E51").Value = "  -0.44%  "
